# Add a paragraph of notes below the "Logging" heading and a new
# "Code Review" Heading 2 paragraph after it.

$d = $word.ActiveDocument

# Locate the existing "Logging" Heading 2 paragraph.
$loggingHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 2" -and `
        $p.Range.Text.TrimEnd([char]13) -eq "Logging") {
        $loggingHeading = $p
        break
    }
}

if ($loggingHeading -eq $null) {
    throw "Could not find the 'Logging' Heading 2 paragraph"
}

# Immediately after the heading sits a pre-existing blank paragraph that must
# stay put; the new content is inserted right after that blank paragraph
# (i.e. before the next one, which we recreate untouched once we're done).
$blankAfterHeading = $loggingHeading.Next()
$insertionTarget = $blankAfterHeading.Next()

$r = $insertionTarget.Range
$r.Collapse(1)

$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$bodyParagraph = '<w:p xmlns:w="' + $w + '">' + `
    '<w:r><w:t xml:space="preserve">Logging is valuable for understanding the events that occur while running your program. For example, if you run your model </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>over night</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> and see that it''s producing ridiculous results the next day, log messages can really help you understand more about the context in which this occurred.</w:t></w:r>' + `
    '</w:p>'

$blankParagraph = '<w:p xmlns:w="' + $w + '"/>'

$codeReviewHeading = '<w:p xmlns:w="' + $w + '">' + `
    '<w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' + `
    '<w:r><w:t>Code Review</w:t></w:r>' + `
    '</w:p>'

# InsertXML on a collapsed range replaces the (empty) target paragraph, so we
# re-emit a trailing blank paragraph to stand in for the one being replaced.
$xml = $bodyParagraph + $blankParagraph + $codeReviewHeading + $blankParagraph

[void]$r.InsertXML($xml)

Write-Output "Added logging notes paragraph and 'Code Review' heading."
